# Auto-generated Excel COM-interop edit script
# Applies the "cryptos" price/volume update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing Excel to keep it as literal text
# (prevents numeric/date auto-conversion of strings such as "1.000" or "26.869.54"),
# then restores the cells original style so no stray formatting is introduced.
function Set-TextValue {
    param($Address, $Text)
    $rng = $ws.Range($Address)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = $origStyle
}

Set-TextValue "D2" "26.869.54"
$ws.Range("E2").Value = "  -1.04%  "
Set-TextValue "D3" "1.806.48"
$ws.Range("E3").Value = "  -0.89%  "
Set-TextValue "D4" "1.000"
$ws.Range("E4").Value = "  -0.59%  "
Set-TextValue "D5" "310.23"
$ws.Range("E5").Value = "  -0.85%  "
Set-TextValue "D6" "1.000"
$ws.Range("E6").Value = "  -0.49%  "
Set-TextValue "D7" "0.4440"
$ws.Range("E7").Value = "  +4.93%  "
Set-TextValue "D8" "0.3677"
$ws.Range("E8").Value = "  -0.42%  "
Set-TextValue "D9" "0.07343"
$ws.Range("E9").Value = "  +1.34%  "
Set-TextValue "D10" "0.8569"
$ws.Range("E10").Value = "  -0.12%  "
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D11" "2.010.87"
$ws.Range("E11").Value = "  +10.24%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D12" "20.66"
$ws.Range("E12").Value = "  -1.47%  "
Set-TextValue "D13" "6.604"
$ws.Range("E13").Value = "  -1.53%  "
Set-TextValue "D14" "92.59"
$ws.Range("E14").Value = "  +3.44%  "
Set-TextValue "D15" "5.303"
$ws.Range("E15").Value = "  +0.00%  "
Set-TextValue "D16" "0.07058"
$ws.Range("E16").Value = "  -0.44%  "
Set-TextValue "D17" "1.001"
$ws.Range("E17").Value = "  -0.60%  "
Set-TextValue "D18" "0.000008732"
$ws.Range("E18").Value = "  -1.33%  "
$ws.Range("E19").Value = "  -0.44%  "
Set-TextValue "D20" "14.87"
$ws.Range("E20").Value = "  -1.05%  "
Set-TextValue "D21" "26.891.92"
$ws.Range("E21").Value = "  -1.24%  "
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("E23").Value = "  -0.86%  "
Set-TextValue "D24" "1.988"
$ws.Range("E24").Value = "  -0.02%  "
Set-TextValue "D25" "151.85"
$ws.Range("E25").Value = "  -0.47%  "
Set-TextValue "D26" "18.50"
$ws.Range("E26").Value = "  +0.58%  "
Set-TextValue "D27" "2.183"
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("E28").Value = "  -0.16%  "
Set-TextValue "D29" "116.56"
$ws.Range("E29").Value = "  +0.18%  "
Set-TextValue "D30" "0.08835"
$ws.Range("E30").Value = "  -0.09%  "
Set-TextValue "D31" "0.7497"
$ws.Range("E31").Value = "  -0.07%  "
Set-TextValue "D32" "1.175"
$ws.Range("E32").Value = "  -1.21%  "
Set-TextValue "D33" "2.931"
$ws.Range("E33").Value = "  +4.57%  "
Set-TextValue "D34" "4.465"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("E35").Value = "  -0.54%  "
$ws.Range("E36").Value = "  -2.90%  "
Set-TextValue "D37" "0.01966"
$ws.Range("E37").Value = "  -0.20%  "
Set-TextValue "D38" "0.05192"
$ws.Range("E38").Value = "  -1.13%  "
Set-TextValue "D39" "0.5314"
$ws.Range("E39").Value = "  +5.55%  "
Set-TextValue "D40" "2.862"
$ws.Range("E40").Value = "  -0.69%  "
Set-TextValue "D41" "7.022"
$ws.Range("E41").Value = "  -3.97%  "
Set-TextValue "D42" "0.1690"
$ws.Range("E42").Value = "  -0.18%  "
Set-TextValue "D43" "0.5158"
$ws.Range("E43").Value = "  +8.88%  "
Set-TextValue "D44" "8.424"
$ws.Range("E44").Value = "  -2.74%  "
Set-TextValue "D45" "1.986"
$ws.Range("E45").Value = "  +7.14%  "
Set-TextValue "D46" "10.57"
$ws.Range("E46").Value = "  -0.47%  "
Set-TextValue "D47" "105.24"
$ws.Range("E47").Value = "  -1.16%  "
Set-TextValue "D48" "0.9997"
$ws.Range("E48").Value = "  -0.56%  "
Set-TextValue "D49" "1.665"
$ws.Range("E49").Value = "  -0.14%  "
Set-TextValue "D50" "0.06319"
$ws.Range("E50").Value = "  -1.12%  "
Set-TextValue "D51" "0.9185"
$ws.Range("E51").Value = "  +0.66%  "
